$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

$ws.Range("B11").Value = 0.08781567715392302
$ws.Range("B12").Value = 0.2243203987552034
$ws.Range("C12").Value = "{'codebleu': 0.22432039875520338, 'ngram_match_score': 0.08764631151099843, 'weighted_ngram_match_score': 0.13212652987798904, 'syntax_match_score': 0.4387990762124711, 'dataflow_match_score': 0.23870967741935484}"
$ws.Range("B13").Value = 0.8943429703917645
